# 2 May 2019 meeting minutes
# Follow-up meeting with COS301 to ascertain current progress, and clarify
# deliverables for tomorrow.
#
# The underlying content edit is: the participant name "Christiaan" was
# changed to "Christo". Because the document's hidden "_GoBack" bookmark
# (which Word drops at the location of the most recent edit) previously sat
# at the very start of the document, this edit relocates it to sit right
# after the newly typed text - which is also where the run ends up getting
# split in two ("Christ" / "o") because the bookmark markers land between
# them while the edit is being made character-by-character.

$d = $word.ActiveDocument

# Remove the existing (stale) _GoBack bookmark - it currently wraps almost
# the whole document (from just before "Summary" to just before the final
# "... heading" paragraph's end).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$text = $d.Content.Text
$idx = $text.IndexOf("Christiaan")
if ($idx -lt 0) {
    throw "Could not locate 'Christiaan' in the document content."
}

# Split "Christiaan" into "Christ" / "iaan" by dropping a collapsed bookmark
# between them - inserting bookmark start/end markers in the middle of a
# run's text forces the run to split in two around the markers.
$splitPoint = $d.Range($idx + 6, $idx + 6)
$d.Bookmarks.Add("_GoBack", $splitPoint)

# Replace the trailing "iaan" with "o", turning "Christiaan" into "Christo"
# while leaving the "Christ" / (new) run boundary - and the bookmark that
# currently sits right before it - untouched.
$tail = $d.Range($idx + 6, $idx + 10)
$tail.Text = "o"

# The bookmark is currently sitting between "Christ" and "o"; it belongs
# after "o" instead. Move it there: drop it, type a throwaway placeholder
# character after "o" so the target position is no longer the very end of
# the paragraph, re-add the (now collapsed) bookmark right before the
# placeholder, then remove the placeholder again.
$d.Bookmarks("_GoBack").Delete()

$afterO = $idx + 7
$placeholder = $d.Range($afterO, $afterO)
$placeholder.InsertAfter("X")

$bmSpot = $d.Range($afterO, $afterO)
$d.Bookmarks.Add("_GoBack", $bmSpot)

$placeholderRange = $d.Range($afterO, $afterO + 1)
$placeholderRange.Delete()
